$d = $word.ActiveDocument

# Title (appears twice: the page Heading1 and the bold SEO-title run near the end).
# ReplaceAll (last arg = 2) takes care of both occurrences in a single pass.
$d.Content.Find.Execute(
    "Play Lightning Leopard Free, Review & Features | Lightning Box Gaming",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Lightning Leopard Slot Game for Free", 2)

# "What we like" bullet: soundtrack -> theme
$d.Content.Find.Execute(
    "Typical soundtrack of Nepal adds to the immersive experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Attractive theme set in Nepal", 2)

# "What we like" bullet: jackpots -> animals/symbols
$d.Content.Find.Execute(
    "Exclusive chance to win jackpots during bonus rounds",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Variety of animals and symbols", 2)

# "What we don't like" bullet: minimum bet -> limited betting range
$d.Content.Find.Execute(
    "High minimum bet of €0.40",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited betting range", 2)

# "What we don't like" bullet: RTP -> lower RTP
$d.Content.Find.Execute(
    "Average RTP of 95.52%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lower RTP compared to some other games", 2)

# Italic meta description run
$d.Content.Find.Execute(
    "Read our review of Lightning Leopard slot game by Lightning Box Gaming. Play free and discover exclusive special features & bonuses, RTP, and jackpots.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Lightning Leopard, a slot game set in Nepal. Play for free and enjoy impressive graphics and exciting features.", 2)
